# Generate Report for Handoff
# - Overview sheet: "Latest HO Xliff Generate Date" (col G) for the affected rows moves
#   from 2016-08-23 14:35:40 to 2016-08-23 14:35:56
# - zh-cn sheet: Priority (col E) goes from "low" to "ht", and
#   "Latest Handoff Datetime" (col H) moves from 2016-08-23 14:35:35 to 2016-08-23 14:35:51
# - de-de sheet: Priority (col E) goes from "low" to "ht", and
#   "Latest Handoff Datetime" (col H) moves from 2016-08-23 14:35:40 to 2016-08-23 14:35:56

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-23 14:35:56"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-23 14:35:51"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-23 14:35:56"
}
